$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats paste-special mode constant (used throughout to replicate a source
# cell's number format / font / alignment onto the new column M cells).
$xlPasteFormats = -4122
$xlRight = -4152

# Row 2: new trailing border-only cell (no value), same format as L2
[void]$ws.Range("L2").Copy()
[void]$ws.Range("M2").PasteSpecial($xlPasteFormats)

# Row 3: M3 = 2021
[void]$ws.Range("L3").Copy()
[void]$ws.Range("M3").PasteSpecial($xlPasteFormats)
$ws.Range("M3").Value = 2021

# Row 4: M4 = 2.017314837395458
[void]$ws.Range("L4").Copy()
[void]$ws.Range("M4").PasteSpecial($xlPasteFormats)
$ws.Range("M4").Value = 2.017314837395458

# Row 5: M5 = 0.11867182493532386
[void]$ws.Range("L5").Copy()
[void]$ws.Range("M5").PasteSpecial($xlPasteFormats)
$ws.Range("M5").Value = 0.11867182493532386

# Row 6: M6 = 3.944091449932318
[void]$ws.Range("L6").Copy()
[void]$ws.Range("M6").PasteSpecial($xlPasteFormats)
$ws.Range("M6").Value = 3.944091449932318

# Row 7: M7 = 0
[void]$ws.Range("L7").Copy()
[void]$ws.Range("M7").PasteSpecial($xlPasteFormats)
$ws.Range("M7").Value = 0

# Row 8: M8 = '-'
[void]$ws.Range("L8").Copy()
[void]$ws.Range("M8").PasteSpecial($xlPasteFormats)
$ws.Range("M8").Value = "-"

# Row 9: M9 = 0
[void]$ws.Range("L9").Copy()
[void]$ws.Range("M9").PasteSpecial($xlPasteFormats)
$ws.Range("M9").Value = 0

# Row 10: M10 = 0.6292103017456653
[void]$ws.Range("L10").Copy()
[void]$ws.Range("M10").PasteSpecial($xlPasteFormats)
$ws.Range("M10").Value = 0.6292103017456653

# Row 11: M11 = '-'
[void]$ws.Range("L11").Copy()
[void]$ws.Range("M11").PasteSpecial($xlPasteFormats)
$ws.Range("M11").Value = "-"

# Row 12: M12 = 1.2497227177719943
[void]$ws.Range("L12").Copy()
[void]$ws.Range("M12").PasteSpecial($xlPasteFormats)
$ws.Range("M12").Value = 1.2497227177719943

# Row 13: M13 = 0.1984453789016842
[void]$ws.Range("L13").Copy()
[void]$ws.Range("M13").PasteSpecial($xlPasteFormats)
$ws.Range("M13").Value = 0.1984453789016842

# Row 14: M14 = '-'
[void]$ws.Range("L14").Copy()
[void]$ws.Range("M14").PasteSpecial($xlPasteFormats)
$ws.Range("M14").Value = "-"

# Row 15: M15 = 0.39861918314956984
[void]$ws.Range("L15").Copy()
[void]$ws.Range("M15").PasteSpecial($xlPasteFormats)
$ws.Range("M15").Value = 0.39861918314956984

# Row 16: M16 = 0
[void]$ws.Range("L16").Copy()
[void]$ws.Range("M16").PasteSpecial($xlPasteFormats)
$ws.Range("M16").Value = 0

# Row 17: M17 = '-'
[void]$ws.Range("L17").Copy()
[void]$ws.Range("M17").PasteSpecial($xlPasteFormats)
$ws.Range("M17").Value = "-"

# Row 18: M18 = 0
[void]$ws.Range("L18").Copy()
[void]$ws.Range("M18").PasteSpecial($xlPasteFormats)
$ws.Range("M18").Value = 0

# Row 19: M19 = 0.8552125203112974
[void]$ws.Range("L19").Copy()
[void]$ws.Range("M19").PasteSpecial($xlPasteFormats)
$ws.Range("M19").Value = 0.8552125203112974

# Row 20: M20 = '-'
[void]$ws.Range("L20").Copy()
[void]$ws.Range("M20").PasteSpecial($xlPasteFormats)
$ws.Range("M20").Value = "-"

# Row 21: M21 = 1.6913581464969858
[void]$ws.Range("L21").Copy()
[void]$ws.Range("M21").PasteSpecial($xlPasteFormats)
$ws.Range("M21").Value = 1.6913581464969858

# Row 22: M22 = 1.8347815875998121
[void]$ws.Range("L22").Copy()
[void]$ws.Range("M22").PasteSpecial($xlPasteFormats)
$ws.Range("M22").Value = 1.8347815875998121

# Row 23: M23 = '-'
[void]$ws.Range("L23").Copy()
[void]$ws.Range("M23").PasteSpecial($xlPasteFormats)
$ws.Range("M23").Value = "-"

# Row 24: M24 = 3.6321107648498847
[void]$ws.Range("L24").Copy()
[void]$ws.Range("M24").PasteSpecial($xlPasteFormats)
$ws.Range("M24").Value = 3.6321107648498847

# Row 25: M25 = 6.121156041530003
[void]$ws.Range("L25").Copy()
[void]$ws.Range("M25").PasteSpecial($xlPasteFormats)
$ws.Range("M25").Value = 6.121156041530003

# Row 26: M26 = '-'
[void]$ws.Range("L26").Copy()
[void]$ws.Range("M26").PasteSpecial($xlPasteFormats)
$ws.Range("M26").Value = "-"
$ws.Range("M26").HorizontalAlignment = $xlRight  # new right-aligned numeric style

# Row 27: M27 = 12.437939862560766
[void]$ws.Range("L27").Copy()
[void]$ws.Range("M27").PasteSpecial($xlPasteFormats)
$ws.Range("M27").Value = 12.437939862560766

# Row 28: M28 = 3.6823562661275693
[void]$ws.Range("L28").Copy()
[void]$ws.Range("M28").PasteSpecial($xlPasteFormats)
$ws.Range("M28").Value = 3.6823562661275693

# Row 29: M29 = 0.6943323387022582
[void]$ws.Range("L29").Copy()
[void]$ws.Range("M29").PasteSpecial($xlPasteFormats)
$ws.Range("M29").Value = 0.6943323387022582

# Row 30: M30 = 7.056499035611798
[void]$ws.Range("L30").Copy()
[void]$ws.Range("M30").PasteSpecial($xlPasteFormats)
$ws.Range("M30").Value = 7.056499035611798

# Row 31: M31 = 2.7447727328177227
[void]$ws.Range("L31").Copy()
[void]$ws.Range("M31").PasteSpecial($xlPasteFormats)
$ws.Range("M31").Value = 2.7447727328177227

# Row 32: M32 = '-'
[void]$ws.Range("L8").Copy()
[void]$ws.Range("M32").PasteSpecial($xlPasteFormats)
$ws.Range("M32").Value = "-"

# Row 33: M33 = 5.641855041937789
[void]$ws.Range("L33").Copy()
[void]$ws.Range("M33").PasteSpecial($xlPasteFormats)
$ws.Range("M33").Value = 5.641855041937789

# Put the saved cursor position back where the source workbook had it
[void]$ws.Range("P6").Select()
